# Auto-generated edit script for cryptos.xlsx update
# Updates coin prices/volume percentages; two pairs of rows swap rank order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$value)
    # Force the written value to remain a text cell (matches the workbook's
    # existing convention of storing every Coin/Link/Price/Volume value as a
    # string) even when the text looks like a pure number (e.g. "1.00"),
    # without permanently altering the cell's style.
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Cells.Item(2,4) '59.470.74'
Set-TextValue $ws.Cells.Item(2,5) '  +0.48%  '
Set-TextValue $ws.Cells.Item(3,4) '2.604.63'
Set-TextValue $ws.Cells.Item(3,5) '  +0.67%  '
Set-TextValue $ws.Cells.Item(4,5) '  +0.25%  '
Set-TextValue $ws.Cells.Item(5,4) '537.91'
Set-TextValue $ws.Cells.Item(5,5) '  +2.85%  '
Set-TextValue $ws.Cells.Item(6,4) '141.44'
Set-TextValue $ws.Cells.Item(6,5) '  +1.52%  '
Set-TextValue $ws.Cells.Item(7,5) '  +0.09%  '
Set-TextValue $ws.Cells.Item(8,5) '  +0.23%  '
Set-TextValue $ws.Cells.Item(9,4) '6.50'
Set-TextValue $ws.Cells.Item(9,5) '  -0.36%  '
Set-TextValue $ws.Cells.Item(10,5) '  +1.02%  '
Set-TextValue $ws.Cells.Item(11,4) '0.334'
Set-TextValue $ws.Cells.Item(11,5) '  +1.38%  '
Set-TextValue $ws.Cells.Item(12,5) '  -1.11%  '
Set-TextValue $ws.Cells.Item(13,4) '3.060.18'
Set-TextValue $ws.Cells.Item(14,4) '59.390.12'
Set-TextValue $ws.Cells.Item(14,5) '  +0.72%  '
Set-TextValue $ws.Cells.Item(15,4) '20.71'
Set-TextValue $ws.Cells.Item(15,5) '  +1.06%  '
Set-TextValue $ws.Cells.Item(16,2) 'ShibaInu'
Set-TextValue $ws.Cells.Item(16,3) 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Cells.Item(16,4) '0.0000133'
Set-TextValue $ws.Cells.Item(16,5) '  +0.46%  '
Set-TextValue $ws.Cells.Item(17,2) 'WrappedEther'
Set-TextValue $ws.Cells.Item(17,3) 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Cells.Item(17,4) '2.583.10'
Set-TextValue $ws.Cells.Item(17,5) '  -0.59%  '
Set-TextValue $ws.Cells.Item(18,4) '340.55'
Set-TextValue $ws.Cells.Item(18,5) '  -0.12%  '
Set-TextValue $ws.Cells.Item(19,4) '4.36'
Set-TextValue $ws.Cells.Item(19,5) '  +1.24%  '
Set-TextValue $ws.Cells.Item(20,4) '10.08'
Set-TextValue $ws.Cells.Item(20,5) '  +0.05%  '
Set-TextValue $ws.Cells.Item(21,5) '  -2.32%  '
Set-TextValue $ws.Cells.Item(22,5) '  +0.01%  '
Set-TextValue $ws.Cells.Item(23,4) '67.14'
Set-TextValue $ws.Cells.Item(23,5) '  +1.19%  '
Set-TextValue $ws.Cells.Item(24,2) 'Kaspa'
Set-TextValue $ws.Cells.Item(24,3) 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Cells.Item(24,4) '0.166'
Set-TextValue $ws.Cells.Item(24,5) '  -1.30%  '
Set-TextValue $ws.Cells.Item(25,2) 'Polygon'
Set-TextValue $ws.Cells.Item(25,3) 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws.Cells.Item(25,4) '0.408'
Set-TextValue $ws.Cells.Item(25,5) '  +0.92%  '
Set-TextValue $ws.Cells.Item(26,5) '  +0.41%  '
Set-TextValue $ws.Cells.Item(27,4) '7.21'
Set-TextValue $ws.Cells.Item(27,5) '  +1.87%  '
Set-TextValue $ws.Cells.Item(28,4) '0.0₃0743'
Set-TextValue $ws.Cells.Item(28,5) '  +2.77%  '
Set-TextValue $ws.Cells.Item(30,4) '1.67'
Set-TextValue $ws.Cells.Item(30,5) '  +5.65%  '
Set-TextValue $ws.Cells.Item(31,4) '5.83'
Set-TextValue $ws.Cells.Item(31,5) '  -1.00%  '
Set-TextValue $ws.Cells.Item(32,4) '18.75'
Set-TextValue $ws.Cells.Item(32,5) '  +0.23%  '
Set-TextValue $ws.Cells.Item(33,4) '150.43'
Set-TextValue $ws.Cells.Item(33,5) '  +0.86%  '
Set-TextValue $ws.Cells.Item(34,4) '3.98'
Set-TextValue $ws.Cells.Item(34,5) '  +0.42%  '
Set-TextValue $ws.Cells.Item(35,5) '  +0.32%  '
Set-TextValue $ws.Cells.Item(36,4) '0.845'
Set-TextValue $ws.Cells.Item(36,5) '  +3.96%  '
Set-TextValue $ws.Cells.Item(37,4) '1.45'
Set-TextValue $ws.Cells.Item(37,5) '  -0.87%  '
Set-TextValue $ws.Cells.Item(38,4) '0.823'
Set-TextValue $ws.Cells.Item(38,5) '  -0.16%  '
Set-TextValue $ws.Cells.Item(39,5) '  +0.28%  '
Set-TextValue $ws.Cells.Item(40,4) '1.00'
Set-TextValue $ws.Cells.Item(40,5) '  +0.10%  '
Set-TextValue $ws.Cells.Item(41,4) '274.76'
Set-TextValue $ws.Cells.Item(41,5) '  +1.10%  '
Set-TextValue $ws.Cells.Item(42,4) '0.601'
Set-TextValue $ws.Cells.Item(42,5) '  -0.19%  '
Set-TextValue $ws.Cells.Item(43,4) '10.71'
Set-TextValue $ws.Cells.Item(43,5) '  -0.63%  '
Set-TextValue $ws.Cells.Item(44,4) '0.0949'
Set-TextValue $ws.Cells.Item(44,5) '  -0.20%  '
Set-TextValue $ws.Cells.Item(45,5) '  +1.39%  '
Set-TextValue $ws.Cells.Item(46,5) '  +0.92%  '
Set-TextValue $ws.Cells.Item(47,2) 'Maker'
Set-TextValue $ws.Cells.Item(47,3) 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Cells.Item(47,4) '1.936.23'
Set-TextValue $ws.Cells.Item(47,5) '  -1.68%  '
Set-TextValue $ws.Cells.Item(48,2) 'InjectiveProtocol'
Set-TextValue $ws.Cells.Item(48,3) 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Cells.Item(48,4) '18.44'
Set-TextValue $ws.Cells.Item(48,5) '  +2.57%  '
Set-TextValue $ws.Cells.Item(49,5) '  -1.41%  '
Set-TextValue $ws.Cells.Item(50,4) '111.26'
Set-TextValue $ws.Cells.Item(50,5) '  -2.21%  '
Set-TextValue $ws.Cells.Item(51,5) '  +1.80%  '
